# Generate Report for Handoff
# Regenerates the localization-status report: the first file (b9be4aa0...)
# is now ready for handoff (status flips from "Handed back" to
# "Ready for handoff"), so its handback columns (Latest Target File /
# Latest Handback File) are no longer populated, and the second file's
# handle (e657794e...) gets reassigned to a brand-new guid
# (ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f).

$wb = $excel.ActiveWorkbook

$newFileA  = "50527308-acf1-477c-ac29-3589133d0d67.md"
$newFileB  = "ffffc185753e-54b2-4141-865f-c6f8dc0b7d5f.md"
$statusTxt = "Ready for handoff"
$overviewDate = "2016-03-21 10:58:35"

$zhXlf = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.zh-cn.xlf"
$deXlf = "50527308-acf1-477c-ac29-3589133d0d67.01ab74c27eb71b052aae46b9eb042c9b346bfa2d.de-de.xlf"
$handoffDt = "2016-03-21 10:58:32"
$handbackDt = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newFileA
$ov.Range("B2").Value = $statusTxt
$ov.Range("C2").Value = $statusTxt
$ov.Range("D2").Value = $overviewDate

$ov.Range("A3").Value = $newFileB
$ov.Range("B3").Value = $statusTxt
$ov.Range("C3").Value = $statusTxt
$ov.Range("D3").Value = $overviewDate

# Hyperlinks on Overview keep the same targets (rId2 / rId3) - only the
# display text needs to track the new file names.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8827c79e062a6dfbe0ccf7be5bfcef1f71e800b5/e2e/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.md", "", "", $newFileA)
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8827c79e062a6dfbe0ccf7be5bfcef1f71e800b5/e2e/e657794e-219e-4342-9713-9de29747f114.md", "", "", $newFileB)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newFileA
$zh.Range("D2").Value = $zhXlf
$zh.Range("E2").Value = $handoffDt
$zh.Range("F2").ClearContents()
$zh.Range("G2").ClearContents()
$zh.Range("H2").Value = $handbackDt

$zh.Range("A3").Value = $newFileB
$zh.Range("D3").Value = $zhXlf
$zh.Range("E3").Value = $handoffDt
$zh.Range("F3").ClearContents()
$zh.Range("G3").ClearContents()
$zh.Range("H3").Value = $handbackDt

# Rebuild the hyperlinks collection: keep the same rId2..rId5 targets
# (the report generator never re-points these), drop the handback-stage
# links that used to live in F/G, and refresh display text everywhere.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8827c79e062a6dfbe0ccf7be5bfcef1f71e800b5/e2e/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.md", "", "", $newFileA)
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10ea8903f9918d668841426579b45b6ce030a9c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.39f176f70bd8afbb1dd150108cbac2dc02e9e7a8.zh-cn.xlf", "", "", $zhXlf)
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6f3dd2281822457dda390ab165896816bbb8a32c/e2e/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.md", "", "", $newFileB)
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2b4d6ab8bbd6f3708889c0dc0819aa8961f64f36/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.39f176f70bd8afbb1dd150108cbac2dc02e9e7a8.zh-cn.xlf", "", "", $zhXlf)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newFileA
$de.Range("D2").Value = $deXlf
$de.Range("E2").Value = $overviewDate
$de.Range("F2").ClearContents()
$de.Range("G2").ClearContents()
$de.Range("H2").Value = $handbackDt

$de.Range("A3").Value = $newFileB
$de.Range("D3").Value = $deXlf
$de.Range("E3").Value = $overviewDate
$de.Range("F3").ClearContents()
$de.Range("G3").ClearContents()
$de.Range("H3").Value = $handbackDt

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8827c79e062a6dfbe0ccf7be5bfcef1f71e800b5/e2e/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.md", "", "", $newFileA)
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a70b8fb572a23b5fe0019416e20467ff9251b9c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.39f176f70bd8afbb1dd150108cbac2dc02e9e7a8.de-de.xlf", "", "", $deXlf)
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/17489a2ac281d3ff98e5294397d0de653fd9e5f2/e2e/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.md", "", "", $newFileB)
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90c86d9c63c0e14d9ccb7691497f52b520b5ca1a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b9be4aa0-ecc7-43c8-bd48-d5d7c8be8f19.39f176f70bd8afbb1dd150108cbac2dc02e9e7a8.de-de.xlf", "", "", $deXlf)

Write-Output "Report regenerated for handoff."
